# Insert a new weekly price record for "Perejil" (Parsley) at row 237 of the
# "Terminal La Palmera de La Serena" sheet. Inserting the row shifts every
# following row down by one (old row 237 -> 238, ..., old row 259 -> 260),
# which matches the rest of the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(237).Insert()

$ws.Cells.Item(237, 1).Value  = 8
$ws.Cells.Item(237, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(237, 3).Value  = "Coquimbo"
$ws.Cells.Item(237, 4).Value  = 45223
$ws.Cells.Item(237, 5).Value  = 4
$ws.Cells.Item(237, 6).Value  = 100112044
$ws.Cells.Item(237, 7).Value  = "Perejil"
$ws.Cells.Item(237, 8).Value  = "Sin especificar"
$ws.Cells.Item(237, 9).Value  = "Primera"
$ws.Cells.Item(237, 10).Value = 2000
$ws.Cells.Item(237, 11).Value = 1400
$ws.Cells.Item(237, 12).Value = 1500
$ws.Cells.Item(237, 13).Value = 1450
$ws.Cells.Item(237, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(237, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(237, 16).Value = 967
$ws.Cells.Item(237, 17).Value = 1.5
$ws.Cells.Item(237, 18).Value = "Hortaliza"
